# Applies the docx edit: extends the final paragraph and appends the
# "Bai 42" lesson block (binding san pham va Slide tren trang chu).

$d = $word.ActiveDocument

# --- Step 1: insert the new text right after "model.content" (same run/paragraph) ---
$rng = $d.Content
$null = $rng.Find.Execute("model.content", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(' > thay đổi thuộc tính bằng placeholder')

# manual line break, then the remainder of the paragraph
$rng.Collapse(0)
$brk = [char]11
$rng.InsertAfter($brk)
$rng.Collapse(0)
$rng.InsertAfter('> thêm vào AutoMapperConfig đối tượng mới > ')

# --- Step 2: drop the now-orphaned trailing manual break run ---
$p172 = $d.Paragraphs.Item(172)
$pr = $p172.Range
$n = $pr.Characters.Count
$delRng = $d.Range($pr.Start + $n - 2, $pr.Start + $n - 1)
$delRng.Delete()

# --- Step 3: append the new paragraphs (all "No Spacing" style, like the lesson headers) ---
$p172 = $d.Paragraphs.Item(172)
$p172.Range.InsertParagraphAfter()

$pBlank1 = $d.Paragraphs.Item(173)
$pBlank1.Style = "No Spacing"

$pBlank1.Range.InsertParagraphAfter()
$pHeading = $d.Paragraphs.Item(174)
$pHeading.Range.Text = 'Bài 42 : Binding sản phẩm và Slide trên trang chủ'
$pHeading.Style = "No Spacing"

$pHeading.Range.InsertParagraphAfter()
$pBlank2 = $d.Paragraphs.Item(175)
$pBlank2.Style = "No Spacing"

$pBlank2.Range.InsertParagraphAfter()
$pBody = $d.Paragraphs.Item(176)
$pBody.Range.Text = 'Thực hành:  vào data Migration dung phương thức Seed để tạo mới CreateSlide >  thêm Prop Content vào Model Side  > add migration addContentToSlides và update Database > Thay đổi layout.html sử dụng thuộc tính RenderSection(“HomeProduct”) > Khai bái Section HomeProduct vào index html > Vì trường hợp một view sử dụng 2 model có thể dung ViewBag, nhưng sẽ làm xấu code nên sẽ tạo HomeViewModel và SlideViewModel > Vào Class AutoMapper cấu hình cho ViewModel mới > tiến hành khai báo Index của HomeController sử dụng HomeViewModel và lặp dữ liệu các model Slide và Product > đặt thêm RenderSection(“FooterJS”) trong layout và tạo SectionFooterJS ở Index Home > Tiêm Silde vào CommonService > Khai báo homeController > viết thêm phương thức cho ProductService > '
$pBody.Style = "No Spacing"

Write-Output ("Done. Paragraph count: " + $d.Paragraphs.Count)
